# dropped NaN's from grade excel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header row: A1 = SIS User ID, B1 = Exams Final Score (same style as before)
$ws.Range("A1").Value = "SIS User ID"
$ws.Range("B1").Value = "Exams Final Score"

# Replace student-name text in column A with their numeric SIS User IDs
$ws.Range("A2").Value = 11
$ws.Range("A3").Value = 12
$ws.Range("A4").Value = 13
$ws.Range("A5").Value = 14
$ws.Range("A6").Value = 15

# Add a new bold header "Projects Final Score" in D1 and widen the column
$ws.Range("D1").Value = "Projects Final Score"
$ws.Range("D1").Font.Bold = $true
$ws.Columns.Item(4).ColumnWidth = 17.65

# Select D1 to match the saved selection state
$ws.Range("D1").Select()
